$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set final cell values (rows 1-24) ---
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "LOM3239"
$ws.Range("C2").Value = "LOM3239"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Projeto Integrado II"
$ws.Range("C3").Value = " Projeto Integrado II"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Integrated Project II"
$ws.Range("C4").Value = "Integrated Project II"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "1"
$ws.Range("C5").Value = "1"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "4"
$ws.Range("C6").Value = "4"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "135 h"
$ws.Range("C7").Value = "135 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2012"
$ws.Range("C8").Value = "01/01/2012"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EF-8"
$ws.Range("C9").Value = "EF-8"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas, reuniões com supervisor, desenvolvimento e elaboração de projeto."
$ws.Range("C19").Value = "Aulas expositivas, reuniões com supervisor, desenvolvimento e elaboração de projeto."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Nota de avaliação do projeto."
$ws.Range("C20").Value = "Nota de avaliação do projeto."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A critério da Comissão de Curso poderá ser oferecida recuperação."
$ws.Range("C21").Value = "A critério da Comissão de Curso poderá ser oferecida recuperação."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOQ4050 -  Engenharia Econômica  (Requisito)`n"
$ws.Range("C23").Value = "LOQ4050 -  Engenharia Econômica  (Requisito)`n"
$ws.Range("B24").Value = "LOQ4234 -  Empreendedorismo  (Requisito)`n"
$ws.Range("C24").Value = "LOQ4234 -  Empreendedorismo  (Requisito)`n"

# --- Clear stale cells left over from the old layout ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Row heights ---
$ws.Rows("10:10").RowHeight = 60
$ws.Rows("11:11").RowHeight = 60
$ws.Rows("13:13").RowHeight = 60
$ws.Rows("14:14").RowHeight = 60
$ws.Rows("15:15").RowHeight = 120
$ws.Rows("16:16").RowHeight = 120
$ws.Rows("18:18").RowHeight = 60
$ws.Rows("19:19").RowHeight = 60
$ws.Rows("20:20").RowHeight = 60
$ws.Rows("21:21").RowHeight = 120
$ws.Rows("23:23").RowHeight = 30
$ws.Rows("24:24").RowHeight = 30

# --- Reset rows that must go back to default (no custom height) ---
$ws.Rows("17:17").AutoFit()
$ws.Rows("22:22").AutoFit()

# --- Remove the now-unused trailing row 25 ---
$ws.Rows("25:25").Delete()
